# Log file updated, with links of Post62
# Adds a new row (row 72) to the blog log table on Sheet1 for "Post62":
#   S.No=62, Title, Date of Post=2020-12-14 (serial 44179),
#   Link on Hashnode, Link on Dev.to

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Duplicate the previous data row (61) into the new row (72) so that the
# new row picks up the same cell formatting (borders / date format /
# hyperlink-like styling) that the rest of the table uses.
$ws.Range("B71:F71").Copy($ws.Range("B72:F72"))

# Now overwrite the copied values with the new Post62 data. The write
# order below matches column order left-to-right except Title/Hashnode
# link swapped, reproducing the shared-string insertion order used by
# the source workbook (Hashnode link, then Title, then Dev.to link).
$ws.Range("B72").Value = 62
$ws.Range("E72").Value = "https://programmingport.hashnode.dev/read-files-content-using-while-loop-or-shell-scripting"
$ws.Range("C72").Value = "Read File's Content using While Loop | Shell Scripting"
$ws.Range("D72").Value = 44179
$ws.Range("F72").Value = "https://dev.to/rahulmishra05/read-file-s-content-using-while-loop-shell-scripting-2anl"

# Grow the "Table2" structured table so it covers the new row too.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("B10:F72"))

# Match the author's selection/active-cell state after the edit.
[void]$ws.Range("F72").Select()
